$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.675.70"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "'3.033.28"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'379.43"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'103.36"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D10").Value = "'36.75"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "'0.0861"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "'3.516.18"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "'18.57"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").Value = "'7.76"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "'3.025.69"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("E18").Value = "  -11.56%  "
$ws.Range("D19").Value = "'51.683.16"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "'12.52"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'0.0₃0963"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "'269.15"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").Value = "'8.22"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").Value = "'7.62"
$ws.Range("E27").Value = "  +6.46%  "
$ws.Range("E28").Value = "  +6.09%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'26.26"
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'10.28"
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").Value = "'34.23"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("D34").Value = "'50.56"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'0.0450"
$ws.Range("E36").Value = "  +3.79%  "
$ws.Range("E38").Value = "  +7.20%  "
$ws.Range("D39").Value = "'0.288"
$ws.Range("E39").Value = "  +9.83%  "
$ws.Range("D40").Value = "'17.16"
$ws.Range("E40").Value = "  +3.29%  "
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").Value = "'2.60"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'127.33"
$ws.Range("E43").Value = "  +7.82%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.116"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("E45").Value = "  +7.35%  "
$ws.Range("D46").Value = "'21.80"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("E47").Value = "  +4.16%  "
$ws.Range("D48").Value = "'2.40"
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("D49").Value = "'2.034.75"
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("D50").Value = "'3.330.73"
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("E51").Value = "  +1.80%  "
